# Applies the cryptos.xlsx price/volume refresh described by the commit
# "Updated cryptos list on Sun Mar 17 10:25:12 UTC 2024 with GitHub Actions".
#
# Column D ("Price") and column E ("Volume(1h)") are stored as literal text
# (not numbers) in the source sheet, e.g. "0.998" or "  -3.18%  ", so every
# write below forces a Text number format before assigning the string and
# then restores the default "Normal" style, which keeps Excel from silently
# re-typing numeric-looking text (like "0.997" or "42.10") as a float and
# from leaving a stray custom style on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# --- Update Price (column D) and Volume(1h) (column E) for the refreshed rows ---
Set-TextValue 2 4 '67.006.70'
Set-TextValue 2 5 '  -2.33%  '

Set-TextValue 3 4 '3.581.03'
Set-TextValue 3 5 '  -3.51%  '

Set-TextValue 4 4 '0.997'
Set-TextValue 4 5 '  -0.31%  '

Set-TextValue 5 4 '575.28'
Set-TextValue 5 5 '  -7.51%  '

Set-TextValue 6 4 '188.93'
Set-TextValue 6 5 '  -3.17%  '

Set-TextValue 7 4 '3.577.75'
Set-TextValue 7 5 '  -3.41%  '

Set-TextValue 8 4 '0.618'
Set-TextValue 8 5 '  -2.94%  '

Set-TextValue 9 4 '0.998'
Set-TextValue 9 5 '  +0.32%  '

Set-TextValue 10 4 '0.678'
Set-TextValue 10 5 '  -7.12%  '

Set-TextValue 11 4 '0.150'
Set-TextValue 11 5 '  -6.49%  '

Set-TextValue 12 4 '55.90'
Set-TextValue 12 5 '  -8.01%  '

Set-TextValue 13 4 '0.0000267'
Set-TextValue 13 5 '  -7.37%  '

Set-TextValue 14 4 '9.87'
Set-TextValue 14 5 '  -5.84%  '

Set-TextValue 15 4 '4.147.99'
Set-TextValue 15 5 '  -3.53%  '

Set-TextValue 16 4 '3.571.20'
Set-TextValue 16 5 '  -3.71%  '

Set-TextValue 17 5 '  -1.58%  '

Set-TextValue 18 4 '18.38'
Set-TextValue 18 5 '  -5.73%  '

Set-TextValue 19 4 '66.851.25'
Set-TextValue 19 5 '  -2.45%  '

Set-TextValue 20 4 '12.20'
Set-TextValue 20 5 '  -5.50%  '

Set-TextValue 21 5 '  -7.63%  '

Set-TextValue 22 4 '396.59'
Set-TextValue 22 5 '  -3.38%  '

Set-TextValue 23 4 '4.21'
Set-TextValue 23 5 '  -10.90%  '

Set-TextValue 24 4 '85.73'
Set-TextValue 24 5 '  -5.03%  '

Set-TextValue 25 4 '11.29'
Set-TextValue 25 5 '  -2.97%  '

Set-TextValue 26 4 '2.93'
Set-TextValue 26 5 '  -5.35%  '

Set-TextValue 27 4 '12.45'
Set-TextValue 27 5 '  -5.06%  '

Set-TextValue 28 4 '6.06'
Set-TextValue 28 5 '  +0.40%  '

Set-TextValue 29 4 '3.62'
Set-TextValue 29 5 '  -4.89%  '

Set-TextValue 30 4 '8.93'
Set-TextValue 30 5 '  -7.84%  '

Set-TextValue 31 4 '7.61'
Set-TextValue 31 5 '  -0.87%  '

Set-TextValue 32 4 '31.04'
Set-TextValue 32 5 '  -5.58%  '

Set-TextValue 33 4 '637.26'
Set-TextValue 33 5 '  +0.10%  '

Set-TextValue 34 4 '12.21'
Set-TextValue 34 5 '  -4.33%  '

Set-TextValue 35 5 '  -7.22%  '

Set-TextValue 36 4 '63.56'
Set-TextValue 36 5 '  -6.49%  '

Set-TextValue 37 4 '42.10'
Set-TextValue 37 5 '  -9.68%  '

Set-TextValue 38 4 '0.403'
Set-TextValue 38 5 '  -2.71%  '

Set-TextValue 39 5 '  +0.27%  '

Set-TextValue 40 4 '0.0₃0759'
Set-TextValue 40 5 '  -8.35%  '


# --- Rows 41/42: Maker and Kaspa swapped rank positions (Coin/Link/Price/Volume) ---
Set-TextValue 41 2 'Kaspa'
Set-TextValue 41 3 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 41 4 '0.133'
Set-TextValue 41 5 '  -5.18%  '

Set-TextValue 42 2 'Maker'
Set-TextValue 42 3 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 42 4 '3.117.73'
Set-TextValue 42 5 '  +6.55%  '

Set-TextValue 43 4 '0.996'
Set-TextValue 43 5 '  -0.41%  '

Set-TextValue 44 4 '2.69'
Set-TextValue 44 5 '  +2.52%  '

Set-TextValue 45 5 '  -3.45%  '

Set-TextValue 46 4 '0.0414'
Set-TextValue 46 5 '  -7.66%  '

Set-TextValue 47 5 '  -6.63%  '

Set-TextValue 48 4 '3.08'
Set-TextValue 48 5 '  +0.61%  '

Set-TextValue 49 4 '140.15'
Set-TextValue 49 5 '  -4.46%  '

Set-TextValue 50 4 '8.51'
Set-TextValue 50 5 '  -10.06%  '

Set-TextValue 51 4 '2.76'
Set-TextValue 51 5 '  -0.35%  '
